$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching style of existing header row (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerStyleSource = $ws.Range("H1")
$newHeaders = $ws.Range("I1:J1")
$newHeaders.Font.Bold = $headerStyleSource.Font.Bold
$newHeaders.HorizontalAlignment = $headerStyleSource.HorizontalAlignment
$newHeaders.VerticalAlignment = $headerStyleSource.VerticalAlignment
$newHeaders.Borders.LineStyle = $headerStyleSource.Borders.LineStyle

# Data rows 2-47 for columns I (I0) and J (IF)
$data = @(
    @(2, 7, 8),
    @(3, 6, 7),
    @(4, 7, 7),
    @(5, 7, 7),
    @(6, 8, 8),
    @(7, 6, 6),
    @(8, 8, 8),
    @(9, 6, 6),
    @(10, 11, 12),
    @(11, 7, 7),
    @(12, 6, 7),
    @(13, 6, 7),
    @(14, 9, 10),
    @(15, 7, 8),
    @(16, 7, 7),
    @(17, 9, 9),
    @(18, 7, 8),
    @(19, 6, 6),
    @(20, 8, 8),
    @(21, 6, 6),
    @(22, 8, 8),
    @(23, 7, 8),
    @(24, 3, 4),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 6, 6),
    @(28, 7, 8),
    @(29, 6, 7),
    @(30, 7, 8),
    @(31, 7, 8),
    @(32, 8, 8),
    @(33, 13, 13),
    @(34, 8, 9),
    @(35, 7, 7),
    @(36, 6, 7),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 8, 9),
    @(41, 4, 4),
    @(42, 6, 6),
    @(43, 3, 5),
    @(44, 4, 4),
    @(45, 8, 9),
    @(46, 4, 4),
    @(47, 6, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
